# Pseudonymizer FAKE_DATA key table: add a "STEUERBARESEINKOMMEN" column
# (taxable income) right after VERMÖGEN, and append a trailing "AMOUNT"
# column, shifting the existing HASEL/HASSH columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert new column O; HASEL (old O) / HASSH (old P) shift right to P / Q.
$ws.Columns.Item(15).Insert()
# Append a brand-new trailing column (becomes R after the shift above).
$ws.Columns.Item(18).Insert()

# Headers
$ws.Cells.Item(1, 15).Value = "STEUERBARESEINKOMMEN"
$ws.Cells.Item(1, 18).Value = "AMOUNT"

# Row 2
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 18).Value = -30

# Row 3
$ws.Cells.Item(3, 15).Value = 20000000
$ws.Cells.Item(3, 18).Value = 50000

# Row 4
$ws.Cells.Item(4, 15).Value = 120003

# Row 6
$ws.Cells.Item(6, 15).Value = 25000
$ws.Cells.Item(6, 18).Value = 0

# Row 7
$ws.Cells.Item(7, 15).Value = 500000
$ws.Cells.Item(7, 18).Value = 3600

# Row 8
$ws.Cells.Item(8, 15).Value = 50000
$ws.Cells.Item(8, 18).Value = 1200

[void]$ws.Range("R9").Select()
